$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.897.01'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.632.11'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.77'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.258'
$ws.Range('E9').Value = '  -3.11%  '
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.863.91'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').Value = '1.636.68'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.564'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.23'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '27.897.30'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '230.37'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('E20').Value = '  -2.51%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.36'
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.95'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.63'
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.38'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').Value = '1.399.76'
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('E36').Value = '  +9.37%  '
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.561'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '66.81'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.81'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '1.773.53'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '87.64'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('E51').Value = '  -0.31%  '
